$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting all existing data down by one row
$ws.Rows.Item(1).Insert()

# Set the header value for the new first row
$ws.Range("A1").Value = "id"
